$wb = $excel.ActiveWorkbook

# 1. Refresh the publication date on the Metadata sheet (IG was re-deployed).
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2021-11-22T14:07:32+00:00"

# 2. Add the new "EGA" (Exome germline analysis) concept as row 5 on the
#    Concepts sheet, right after the existing SA / VC / QC rows.
$concepts = $wb.Worksheets.Item("Concepts")

# Format the Level cell as text first so "1" is written as a shared string
# (matching how the other rows store their Level value) instead of a number.
$concepts.Range("A5").NumberFormat = "@"
$concepts.Range("A5").Value = "1"
$concepts.Range("B5").Value = "EGA"
$concepts.Range("C5").Value = "Exome germline analysis"
$concepts.Range("D5").Value = "Exome germline analysis"

# Copy the formatting (borders/fill/alignment) from the row above so the new
# row matches the existing table styling.
$concepts.Range("A4:D4").Copy()
$concepts.Range("A5:D5").PasteSpecial(-4122)
